$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.06127659574468085
$ws.Range("A3").Value = 0.1931914893617021
$ws.Range("A4").Value = 0.4451063829787234
$ws.Range("A5").Value = 0.796595744680851
$ws.Range("A6").Value = 0.9570212765957447
$ws.Range("A7").Value = 0.982127659574468
$ws.Range("A8").Value = 0.9868085106382978
$ws.Range("A9").Value = 0.9957446808510637
$ws.Range("A10").Value = 0.9974468085106382
$ws.Range("A11").Value = 0.9982978723404254
$ws.Range("A12").Value = 0.9982978723404254
$ws.Range("A13").Value = 0.998723404255319
$ws.Range("A14").Value = 0.9991489361702126
$ws.Range("A15").Value = 0.9991489361702126
$ws.Range("A16").Value = 0.9991489361702126
$ws.Range("A17").Value = 0.9991489361702126
$ws.Range("A18").Value = 0.9991489361702126
$ws.Range("A19").Value = 0.9991489361702126
$ws.Range("A20").Value = 0.9991489361702126
$ws.Range("A21").Value = 0.9991489361702126
$ws.Range("A22").Value = 0.9991489361702126
$ws.Range("A23").Value = 0.9991489361702126
$ws.Range("A24").Value = 0.9991489361702126
$ws.Range("A25").Value = 0.9991489361702126
$ws.Range("A26").Value = 0.9991489361702126
$ws.Range("A27").Value = 0.9991489361702126
$ws.Range("A28").Value = 0.9991489361702126
$ws.Range("A29").Value = 0.9991489361702126
$ws.Range("A30").Value = 0.9995744680851062
$ws.Range("A31").Value = 0.9999999999999998
